$p = $ppt.ActivePresentation

# Slide 3 ("Konfigurasi .env"): remove a leftover empty decorative
# rectangle (no visible text, a stray highlight bar) — shape id 65314965,
# the 6th shape on the slide.
$s3 = $p.Slides.Item(3)
for ($i = $s3.Shapes.Count; $i -ge 1; $i--) {
    $sh = $s3.Shapes.Item($i)
    if ($sh.Id -eq 65314965) {
        $sh.Delete()
    }
}

# Slide 4 ("Test Koneksi"): tiny rotation re-normalization of the
# decorative circle picture (cosmetic resave artifact).
$s4 = $p.Slides.Item(4)
for ($i = 1; $i -le $s4.Shapes.Count; $i++) {
    $sh = $s4.Shapes.Item($i)
    if ($sh.Id -eq 923076003) {
        $sh.Rotation = $sh.Rotation
    }
}

# Delete slide 5 ("Composer") entirely — it is a leftover/duplicate
# slide; removing it shifts the following slide ("Laravel 11" /
# "Migrasi") up to become the new slide 5.
$p.Slides.Item(5).Delete()
